# Apply updated cryptos list values (price & 1h volume change) per row.
# Column D (Price) values that look numeric are prefixed with a literal
# leading apostrophe so Excel stores them as text (matching the source
# workbook, where these cells are inline/shared strings, not numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '''61.092.26'
$ws.Range("E2").Value = '  +3.04%  '

# Row 3: Ethereum
$ws.Range("D3").Value = '''3.259.11'
$ws.Range("E3").Value = '  +2.48%  '

# Row 4: TetherUSD
$ws.Range("E4").Value = '  +0.00%  '

# Row 5: BNB
$ws.Range("D5").Value = '''545.88'
$ws.Range("E5").Value = '  +2.41%  '

# Row 6: Solana
$ws.Range("D6").Value = '''148.43'
$ws.Range("E6").Value = '  +4.30%  '

# Row 7: USDC
$ws.Range("E7").Value = '  -0.15%  '

# Row 8: XRP
$ws.Range("E8").Value = '  -0.88%  '

# Row 9: Toncoin
$ws.Range("D9").Value = '''7.44'
$ws.Range("E9").Value = '  +2.30%  '

# Row 10: Dogecoin
$ws.Range("E10").Value = '  +3.00%  '

# Row 11: Cardano
$ws.Range("E11").Value = '  -1.74%  '

# Row 12: WrappedliquidstakedEther2.0
$ws.Range("D12").Value = '''3.822.32'
$ws.Range("E12").Value = '  +2.48%  '

# Row 13: TRON
$ws.Range("E13").Value = '  -1.17%  '

# Row 14: Avalanche
$ws.Range("D14").Value = '''26.56'
$ws.Range("E14").Value = '  +2.32%  '

# Row 15: ShibaInu
$ws.Range("E15").Value = '  +3.00%  '

# Row 16: WrappedBTC
$ws.Range("D16").Value = '''61.091.62'
$ws.Range("E16").Value = '  +3.00%  '

# Row 17: WrappedEther
$ws.Range("D17").Value = '''3.266.65'
$ws.Range("E17").Value = '  +1.79%  '

# Row 18: Polkadot
$ws.Range("D18").Value = '''6.35'
$ws.Range("E18").Value = '  +1.62%  '

# Row 19: Chainlink
$ws.Range("D19").Value = '''13.49'
$ws.Range("E19").Value = '  +3.45%  '

# Row 20: Uniswap
$ws.Range("D20").Value = '''8.46'
$ws.Range("E20").Value = '  +3.34%  '

# Row 21: BitcoinCash
$ws.Range("D21").Value = '''378.16'
$ws.Range("E21").Value = '  +0.30%  '

# Row 22: Dai
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  -0.09%  '

# Row 23: Polygon
$ws.Range("D23").Value = '''0.534'
$ws.Range("E23").Value = '  +0.28%  '

# Row 24: Litecoin
$ws.Range("D24").Value = '''70.16'
$ws.Range("E24").Value = '  +0.48%  '

# Row 25: Kaspa
$ws.Range("D25").Value = '''0.172'
$ws.Range("E25").Value = '  +2.14%  '

# Row 26: InternetComputer(DFINITY)
$ws.Range("D26").Value = '''8.62'
$ws.Range("E26").Value = '  +2.29%  '

# Row 27: Binance-PegBSC-USD
$ws.Range("E27").Value = '  -0.07%  '

# Row 28: PEPE
$ws.Range("D28").Value = '''0.0₃0929'
$ws.Range("E28").Value = '  +6.28%  '

# Row 29: PancakeSwap
$ws.Range("E29").Value = '  +2.54%  '

# Row 30: EthereumClassic
$ws.Range("D30").Value = '''22.60'
$ws.Range("E30").Value = '  +0.60%  '

# Row 31: RenderToken
$ws.Range("D31").Value = '''6.23'
$ws.Range("E31").Value = '  +2.78%  '

# Row 32: NEARProtocol
$ws.Range("D32").Value = '''5.43'
$ws.Range("E32").Value = '  +3.58%  '

# Row 33: Fetch.AI
$ws.Range("D33").Value = '''1.26'
$ws.Range("E33").Value = '  +7.43%  '

# Row 34: Aptos
$ws.Range("D34").Value = '''6.65'
$ws.Range("E34").Value = '  +4.61%  '

# Row 35: Monero
$ws.Range("D35").Value = '''159.37'
$ws.Range("E35").Value = '  +1.29%  '

# Row 36: ImmutableX
$ws.Range("D36").Value = '''1.44'
$ws.Range("E36").Value = '  +7.61%  '

# Row 37: EnergySwap
$ws.Range("D37").Value = '''26.52'
$ws.Range("E37").Value = '  +4.24%  '

# Row 38: Maker
$ws.Range("D38").Value = '''2.806.69'
$ws.Range("E38").Value = '  +3.51%  '

# Row 39: Hedera
$ws.Range("E39").Value = '  +0.81%  '

# Row 40: Stacks
$ws.Range("D40").Value = '''1.74'
$ws.Range("E40").Value = '  +2.34%  '

# Row 41: VeChain
$ws.Range("D41").Value = '''0.0312'
$ws.Range("E41").Value = '  +6.39%  '

# Row 42: Filecoin
$ws.Range("E42").Value = '  +0.32%  '

# Row 43: OKB
$ws.Range("D43").Value = '''40.14'
$ws.Range("E43").Value = '  +2.51%  '

# Row 44: Mantle
$ws.Range("D44").Value = '''0.734'
$ws.Range("E44").Value = '  +1.14%  '

# Row 45: RenzoRestakedETH
$ws.Range("D45").Value = '''3.304.04'
$ws.Range("E45").Value = '  +2.56%  '

# Row 46: Stellar/ONDO (rows swapped)
$ws.Range("B46").Value = 'ONDO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D46").Value = '''1.02'
$ws.Range("E46").Value = '  +2.95%  '

# Row 47: ONDO/Stellar (rows swapped)
$ws.Range("B47").Value = 'Stellar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D47").Value = '''0.105'
$ws.Range("E47").Value = '  +2.06%  '

# Row 48: InjectiveProtocol
$ws.Range("E48").Value = '  +6.42%  '

# Row 49: Cosmos
$ws.Range("D49").Value = '''6.25'
$ws.Range("E49").Value = '  +0.56%  '

# Row 50: SuiNetwork
$ws.Range("D50").Value = '''0.806'
$ws.Range("E50").Value = '  +5.93%  '

# Row 51: Bittensor
$ws.Range("D51").Value = '''279.24'
$ws.Range("E51").Value = '  +8.77%  '
